# LOB1055.xlsx update
#  - Objetivos: fill in the real objective text (was showing the docente name by mistake)
#  - Insert a new row for "Docentes responsáveis:" value (own row under the label),
#    shifting Programa resumido / Short syllabus / Programa / Syllabus / Avaliação /
#    Método / Critério / Norma de recuperação down by one row
#  - Fill in the real "Programa resumido" and "Programa" texts (were showing stray
#    placeholder values)
#  - Re-align Método / Critério / Norma de recuperação with their correct values
#  - Add a new "Bibliografia:" value row at the end

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Objetivos (row 10) — correct text
$ws.Range("B10:C10").Value = "Conscientizar os alunos da importância de uma política de gestão empresarial para assegurar a prevenção de acidentes e doenças do trabalho."

# 2) Insert a blank row at 13 for the "Docentes responsáveis:" value — this shifts the
#    old rows 13-21 down to 14-22, carrying their formatting (styles / row heights) along.
$ws.Rows.Item(13).Insert()

# The freshly-inserted row only carries column A's style by default, so line up B13:C13
# with the normal body formatting (same as B14:C14 / column B-C styles) before filling it in.
$ws.Range("B14:C14").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3) Docentes responsáveis value now lives on its own row (13), under the row-12 label
$ws.Range("B13:C13").Value = "8767640 - Eduardo Ferro dos Santos"

# 4) Programa resumido (row 14) — correct text
$ws.Range("B14:C14").Value = "1) Introdução a Segurança do Trabalho`n2) Programas de Gerenciamento de Riscos`n3) Metodologias de Identificação e Avaliação de Riscos.`n4) Estratégias de Prevenção e Controle de Riscos.`n5) Organização de Serviços de Segurança do Trabalho`n6) Estudos de casos."

# 5) Programa (row 16) — correct (detailed) text
$ws.Range("B16:C16").Value = "1) Introdução a Segurança do Trabalho: Conceitos de acidentes de trabalho, doenças e do papel do engenheiro na segurança do trabalho.`n2) Programas de Gerenciamento de Riscos: Sistemas de gestão em segurança do trabalho, comprometimento e participação gerencial, formalização dos programas, participação dos funcionários, avaliação de programas, coleta e análise de informações.`n3) Metodologias de Identificação e Avaliação de Riscos: Identificação e avaliação de riscos: identificação do risco nas instalações, métodos e técnicas, avaliação da exposição dos funcionários, capacitação ocupacional, acompanhamento de programas e instrumentação aplicada.`n4) Estratégias de Prevenção e Controle de Riscos Revisão do projeto ou modificação nas instalações e nos processos, procedimentos para atividades perigosas, equipamentos de proteção à segurança e à saúde, programa de manutenção preventiva e de ordem e limpeza, investigação de acidentes, segurança das dependências e das áreas restritas, emergências.`n5) Organização de Serviços e da Segurança do Trabalho: Normas regulamentadoras.`n6) Estudos de casos: Aplicações Práticas."

# 6) Método (row 19), Critério (row 20), Norma de recuperação (row 21) — shift each
#    value up to sit under the correct label
$ws.Range("B19:C19").Value = "Duas Notas  N1  1º bimestre e N2  2º bimestre. A composição das N fica a critério do docente."
$ws.Range("B20:C20").Value = "MF = (N1+ N2)/2"
$ws.Range("B21:C21").Value = "NF = (MF + PR)/ 2 , onde PR é uma prova de recuperação"

# 7) Bibliografia (row 22) — new value
$ws.Range("B22:C22").Value = "ABIQUIM . Comissões Técnicas; Guia de implantação saúde e segurança do trabalhador. São Paulo: 1994. 114p.`nFUNDACENTRO. Curso de engenharia de segurança do trabalho. São Paulo, FUNDACENTRO: 1981. 6v. il.`nFUNDACENTRO. Guia para rotulagem preventiva de produtos químicos perigosos. São Pulo, FUNDACENTRO: 1980. 76p.`nHirata, M. H. at all; Manual de biossegurança. Ed. Manole Barueri-SP: 2008. 496p.`nNormas regulamentadoras.< http://portal.mte.gov.br/legislacao/normas-regulamentadoras-1.htm>`nPerez, R. C.; Emergências tecnológicas. Crearte Ed. - Sorocaba SP: 2008. 300p.`nPaoleschi, B.; CIPA: guia prático de segurança do trabalho. Ed. Érica  São Paulo: 2009. 128p.`nSoares, R. A.; Manutenção preventiva. CNI: 1980. 59p.`nSax, N. I. ; Lewis, R. J. ; Dangerous properties of industrial materials. Van Nostrand Reinhold  USA:1989. 3527p.`nMattos, A. de O. (orgs.)- Higiene e segurança do trabalho. Elsevier Ed.  R.J.:2011. 408p. `nVincoli, J. W.; Risk management for hazardous chemicals. CRC Press  USA:1997. 3040P."

Write-Output "LOB1055 updated"
